# Update countries & provincias Spain
# Applies the 9 April 2020 18:22 data refresh to the "Pais" sheet:
#  - Updates the "Datos actualizados..." timestamp in A1
#  - Refreshes numeric figures (Casos totales/Nuevos casos/Casos activos/
#    Recuperados/Casos criticos/Muertes hoy/Muertes) for several countries
#  - Ecuador now overtakes Japon and Pakistan in the ranking (new data,
#    pushing Japon to row 34 and Pakistan to row 35 with their own,
#    unchanged figures)
#  - Argelia now overtakes Islandia in the ranking (new data, pushing
#    Islandia to row 57 with its own, unchanged figures)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Updated timestamp
$ws.Range("A1").Value = "Datos actualizados a 9 de Abril de 2020 a las 18:22"

# row, label, Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes
$rows = @(
    @(4,  "Estados Unidos", 440572, 5645, 24504, 400249, 9318, 1031, 15819),
    @(6,  "Italia",         143626, 4204, 28470,  96877, 3605,  610, 18279),
    @(19, "Austria",         13219,  277,  5240,   7684,  266,   22,   295),
    @(25, "Noruega",          6160,  118,    32,   6020,   78,    7,   108),
    @(33, "Ecuador",          4965,  515,   339,   4354,  139,   30,   272),
    @(34, "Japon",            4667,    0,   632,   3941,   99,    0,    94),
    @(35, "Pakistan",         4489,  226,   572,   3854,   31,    2,    63),
    @(55, "Argentina",        1795,    0,   365,   1359,   96,    6,    71),
    @(56, "Argelia",          1666,   94,   347,   1084,   46,   30,   235),
    @(57, "Islandia",         1648,   32,   688,    954,   11,    0,     6),
    @(81, "Bulgaria",          618,   25,    48,    546,   32,    0,    24)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
    $ws.Cells.Item($r, 8).Value = $row[8]
}
